# envio semana 20 de 2025
# Adds week 20 ("W") column to the weekly IRA-hospital report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new week, matching the look of the other
# week-number headers in row 1 (bold, centered) and stored as text
# (the existing week headers 1..19 are text, not numbers) rather than
# a number.
$ws.Range("W1").Font.Bold = $true
$ws.Range("W1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("W1").Value = "'20"

# Per-facility counts for week 20. Rows that don't carry data for this
# week (rows 3, 15, 16, 18, 19, 25, 31, 40 -- the rows that were already
# missing values in other columns) are intentionally left untouched.
$weekData = @{
    2  = 0
    4  = 0
    5  = 0
    6  = 32
    7  = 2
    8  = 22
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    17 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    26 = 0
    27 = 4
    28 = 10
    29 = 0
    30 = 0
    32 = 16
    33 = 2
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 0
}

foreach ($row in $weekData.Keys) {
    $ws.Cells.Item($row, 23).Value = $weekData[$row]
}
